$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.834.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "'2.328.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'521.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.40%  "
$ws.Range("D6").Value = "'134.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.49%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").Value = "'2.355.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.24%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "'23.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "'2.745.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "'56.877.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "'2.330.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "'10.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "'323.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.03%  "
$ws.Range("D22").Value = "'6.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'60.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").Value = "'0.164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.29%  "
$ws.Range("D26").Value = "'0.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'7.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.78%  "
$ws.Range("E28").Value = "  +14.75%  "
$ws.Range("D29").Value = "'0.0₃0742"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.57%  "
$ws.Range("E30").Value = "  +5.23%  "
$ws.Range("D31").Value = "'166.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").Value = "'6.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'18.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "'1.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").Value = "'0.929"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'4.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.51%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.61%  "
$ws.Range("D40").Value = "'37.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").Value = "'0.384"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").Value = "'3.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.81%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'138.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "'280.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.45%  "
$ws.Range("D46").Value = "'0.0935"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").Value = "'0.0506"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'0.566"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0218"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").Value = "'0.382"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'17.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.80%  "
